$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.293.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.977.51'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.75%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '357.03'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.50'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.91%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.58%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.629'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.11'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0877'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.33%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.43'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.71%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.456.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.10%  '

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.79'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.083.09'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.991'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.311.28'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0987'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.57'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.31'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.81'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.179'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.108'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.47'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.54'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.72%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.17'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0443'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.26'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.14'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.73'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.120'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.61'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.63'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.48'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.55%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.148.96'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0351'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.244'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.928'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.68%  '
